$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving exact formatting (e.g.
# trailing zeros, thousand-dot separators) without leaving the cell's
# NumberFormat/style altered afterwards. We briefly mark the cell as
# Text so Excel doesn't auto-convert number-looking strings (e.g.
# "1.000", "0.9996") into numeric values, then clear the formatting we
# just applied so the cell ends up stateless again, matching the
# original (unstyled) cells.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "29.678.76"

Set-TextValue $ws.Range("D3") "1.854.46"
Set-TextValue $ws.Range("E3") "  +1.11%  "

Set-TextValue $ws.Range("D4") "0.9996"
Set-TextValue $ws.Range("E4") "  +0.05%  "

Set-TextValue $ws.Range("D5") "244.41"
Set-TextValue $ws.Range("E5") "  +0.52%  "

Set-TextValue $ws.Range("D6") "0.6398"
Set-TextValue $ws.Range("E6") "  +3.16%  "

Set-TextValue $ws.Range("E7") "  -0.01%  "

Set-TextValue $ws.Range("E8") "  +3.83%  "

Set-TextValue $ws.Range("D9") "0.07488"
Set-TextValue $ws.Range("E9") "  +1.56%  "

Set-TextValue $ws.Range("D10") "0.2967"
Set-TextValue $ws.Range("E10") "  +2.39%  "

Set-TextValue $ws.Range("D11") "24.33"
Set-TextValue $ws.Range("E11") "  +4.37%  "

Set-TextValue $ws.Range("D12") "0.07653"
Set-TextValue $ws.Range("E12") "  -0.31%  "

Set-TextValue $ws.Range("D13") "1.856.84"
Set-TextValue $ws.Range("E13") "  +1.44%  "

Set-TextValue $ws.Range("D14") "5.029"
Set-TextValue $ws.Range("E14") "  +0.95%  "

Set-TextValue $ws.Range("D15") "0.6889"
Set-TextValue $ws.Range("E15") "  +2.68%  "

Set-TextValue $ws.Range("D16") "83.65"
Set-TextValue $ws.Range("E16") "  +1.17%  "

Set-TextValue $ws.Range("D17") "0.000009680"
Set-TextValue $ws.Range("E17") "  +7.94%  "

Set-TextValue $ws.Range("D18") "6.056"
Set-TextValue $ws.Range("E18") "  +3.08%  "

Set-TextValue $ws.Range("D19") "29.725.27"
Set-TextValue $ws.Range("E19") "  +1.81%  "

Set-TextValue $ws.Range("D20") "2.111.32"
Set-TextValue $ws.Range("E20") "  +2.02%  "

Set-TextValue $ws.Range("D21") "235.58"
Set-TextValue $ws.Range("E21") "  -0.38%  "

Set-TextValue $ws.Range("D22") "12.64"
Set-TextValue $ws.Range("E22") "  +0.99%  "

Set-TextValue $ws.Range("E23") "  +0.01%  "

Set-TextValue $ws.Range("D24") "7.445"
Set-TextValue $ws.Range("E24") "  +1.24%  "

Set-TextValue $ws.Range("D25") "1.000"
Set-TextValue $ws.Range("E25") "  +0.09%  "

Set-TextValue $ws.Range("D26") "158.26"
Set-TextValue $ws.Range("E26") "  -0.01%  "

Set-TextValue $ws.Range("E27") "  +0.80%  "

Set-TextValue $ws.Range("E28") "  -0.59%  "

Set-TextValue $ws.Range("D30") "0.06219"
Set-TextValue $ws.Range("E30") "  +7.86%  "

Set-TextValue $ws.Range("D31") "1.494"
Set-TextValue $ws.Range("E31") "  +0.13%  "

Set-TextValue $ws.Range("D32") "1.277"
Set-TextValue $ws.Range("E32") "  +5.65%  "

Set-TextValue $ws.Range("D33") "4.154"
Set-TextValue $ws.Range("E33") "  +1.50%  "

Set-TextValue $ws.Range("D34") "4.089"
Set-TextValue $ws.Range("E34") "  -0.53%  "

Set-TextValue $ws.Range("E35") "  +1.27%  "

Set-TextValue $ws.Range("E36") "  +2.32%  "

Set-TextValue $ws.Range("D37") "0.7268"
Set-TextValue $ws.Range("E37") "  -0.77%  "

Set-TextValue $ws.Range("E38") "  +0.17%  "

Set-TextValue $ws.Range("E39") "  -1.04%  "

Set-TextValue $ws.Range("D40") "0.01782"
Set-TextValue $ws.Range("E40") "  +1.57%  "

Set-TextValue $ws.Range("D41") "1.200.57"
Set-TextValue $ws.Range("E41") "  -2.15%  "

Set-TextValue $ws.Range("D42") "0.9199"
Set-TextValue $ws.Range("E42") "  +0.82%  "

Set-TextValue $ws.Range("D43") "6.160"
Set-TextValue $ws.Range("E43") "  -1.61%  "

Set-TextValue $ws.Range("E44") "  -0.01%  "

Set-TextValue $ws.Range("D45") "2.016.72"
Set-TextValue $ws.Range("E45") "  +2.19%  "

Set-TextValue $ws.Range("D46") "101.98"
Set-TextValue $ws.Range("E46") "  +0.32%  "

Set-TextValue $ws.Range("D47") "66.38"
Set-TextValue $ws.Range("E47") "  +1.11%  "

Set-TextValue $ws.Range("E48") "  +0.83%  "

Set-TextValue $ws.Range("D49") "9.193"
Set-TextValue $ws.Range("E49") "  +0.18%  "

Set-TextValue $ws.Range("D50") "0.4050"
Set-TextValue $ws.Range("E50") "  +0.49%  "

Set-TextValue $ws.Range("D51") "0.05801"
Set-TextValue $ws.Range("E51") "  +0.97%  "
